# Apply attendance-count updates (0 -> 1) for specific cells on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cell addresses that need to be flipped from 0 to 1.
$cells = @(
    "G3", "H3",
    "D4", "E4",
    "D5", "E5",
    "D6", "E6",
    "H7",
    "H8",
    "D9", "E9",
    "D10", "E10",
    "D11", "E11",
    "D12", "E12",
    "H13",
    "H14",
    "D15", "E15",
    "H16",
    "D17", "E17",
    "H18"
)

foreach ($addr in $cells) {
    $ws.Range($addr).Value = 1
}
